$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 15:05"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1386514
$ws.Range("C4").Value = 680
$ws.Range("E4").Value = 1042437
$ws.Range("G4").Value = 57
$ws.Range("H4").Value = 81852

# --- India (row 15) ---
$ws.Range("D15").Value = 23059
$ws.Range("E15").Value = 46072

# --- Arabia Saudita (row 20) ---
$ws.Range("B20").Value = 42925
$ws.Range("C20").Value = 1911
$ws.Range("D20").Value = 15257
$ws.Range("E20").Value = 27404
$ws.Range("F20").Value = 147
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 264

# --- Argentina (row 56) ---
$ws.Range("E56").Value = 4124
$ws.Range("F56").Value = 170
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 317

# --- Kazajistan (row 60) ---
$ws.Range("F60").Value = 31

# --- Reorder Kenia above Mali (rows 114-117), with refreshed stats ---
# Row 114 becomes Kenia (new, higher counts); Mali/Crucero/Uruguay each
# shift down one row, keeping their own (unchanged) stats.
$ws.Range("A114").Value = "Kenia"
$ws.Range("B114").Value = 715
$ws.Range("C114").Value = 15
$ws.Range("D114").Value = 259
$ws.Range("E114").Value = 420
$ws.Range("F114").Value = 1
$ws.Range("G114").Value = 3
$ws.Range("H114").Value = 36

$ws.Range("A115").Value = "Mali"
$ws.Range("D115").Value = 377
$ws.Range("E115").Value = 296
$ws.Range("F115").Value = 0
$ws.Range("H115").Value = 39

$ws.Range("A116").Value = "Crucero"
$ws.Range("B116").Value = 712
$ws.Range("D116").Value = 651
$ws.Range("E116").Value = 48
$ws.Range("F116").Value = 4
$ws.Range("H116").Value = 13

$ws.Range("A117").Value = "Uruguay"
$ws.Range("B117").Value = 711
$ws.Range("D117").Value = 523
$ws.Range("E117").Value = 169
$ws.Range("F117").Value = 8
$ws.Range("H117").Value = 19
